$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-16 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-17 Tuesday", 2)
$d.Content.Find.Execute("247×9=2223", $true, $false, $false, $false, $false, $true, 1, $false, "780×9=7020", 2)
$d.Content.Find.Execute("574×5=2870", $true, $false, $false, $false, $false, $true, 1, $false, "873×6=5238", 2)
$d.Content.Find.Execute("122×9=1098", $true, $false, $false, $false, $false, $true, 1, $false, "502×9=4518", 2)
$d.Content.Find.Execute("695×4=2780", $true, $false, $false, $false, $false, $true, 1, $false, "938×3=2814", 2)
$d.Content.Find.Execute("833×8=6664", $true, $false, $false, $false, $false, $true, 1, $false, "613×9=5517", 2)
$d.Content.Find.Execute("652×8=5216", $true, $false, $false, $false, $false, $true, 1, $false, "834×6=5004", 2)
$d.Content.Find.Execute("803×7=5621", $true, $false, $false, $false, $false, $true, 1, $false, "333×9=2997", 2)
$d.Content.Find.Execute("992×9=8928", $true, $false, $false, $false, $false, $true, 1, $false, "566×9=5094", 2)
$d.Content.Find.Execute("971×6=5826", $true, $false, $false, $false, $false, $true, 1, $false, "936×6=5616", 2)
$d.Content.Find.Execute("847×5=4235", $true, $false, $false, $false, $false, $true, 1, $false, "439×5=2195", 2)
$d.Content.Find.Execute("143×3=429", $true, $false, $false, $false, $false, $true, 1, $false, "962×8=7696", 2)
$d.Content.Find.Execute("850×6=5100", $true, $false, $false, $false, $false, $true, 1, $false, "904×8=7232", 2)
$d.Content.Find.Execute("861×5=4305", $true, $false, $false, $false, $false, $true, 1, $false, "840×7=5880", 2)
$d.Content.Find.Execute("918×4=3672", $true, $false, $false, $false, $false, $true, 1, $false, "995×6=5970", 2)
$d.Content.Find.Execute("149×8=1192", $true, $false, $false, $false, $false, $true, 1, $false, "505×6=3030", 2)
$d.Content.Find.Execute("237×4=948", $true, $false, $false, $false, $false, $true, 1, $false, "962×9=8658", 2)
$d.Content.Find.Execute("961×3=2883", $true, $false, $false, $false, $false, $true, 1, $false, "746×3=2238", 2)
$d.Content.Find.Execute("460×3=1380", $true, $false, $false, $false, $false, $true, 1, $false, "877×7=6139", 2)
$d.Content.Find.Execute("906×3=2718", $true, $false, $false, $false, $false, $true, 1, $false, "556×4=2224", 2)
$d.Content.Find.Execute("844×8=6752", $true, $false, $false, $false, $false, $true, 1, $false, "475×8=3800", 2)
$d.Content.Find.Execute("107×2=214", $true, $false, $false, $false, $false, $true, 1, $false, "211×3=633", 2)
$d.Content.Find.Execute("570×6=3420", $true, $false, $false, $false, $false, $true, 1, $false, "618×4=2472", 2)
$d.Content.Find.Execute("281×8=2248", $true, $false, $false, $false, $false, $true, 1, $false, "756×6=4536", 2)
$d.Content.Find.Execute("838×2=1676", $true, $false, $false, $false, $false, $true, 1, $false, "853×3=2559", 2)
$d.Content.Find.Execute("890×5=4450", $true, $false, $false, $false, $false, $true, 1, $false, "816×3=2448", 2)
